$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 'Nathanael Wang'
$ws.Range("D3").Value = 'Grace Park'
$ws.Range("E3").Value = 'Michael Dong '
$ws.Range("G3").Value = 'Josh Yang'
$ws.Range("K3").Value = 'chloe lim'
$ws.Range("M3").Value = 'Samuel Wen'
$ws.Range("N3").Value = 'Sabrina Sun'
$ws.Range("C4").Value = 'Christina Ko'
$ws.Range("D4").Value = 'Phillip Seo'
$ws.Range("F4").Value = 'Daniel Kuo'
$ws.Range("G4").Value = 'Susanna Tang'
$ws.Range("K4").Value = 'Sehyun Jung'
$ws.Range("L4").Value = 'Daniel Kim '
$ws.Range("M4").Value = 'Faith Chen'
$ws.Range("N4").Value = 'Josh Yang'
$ws.Range("C5").Value = 'Sehyun Jung'
$ws.Range("D5").Value = 'chloe lim'
$ws.Range("F5").Value = 'Kyle Hwang'
$ws.Range("G5").Value = 'Sabrina Sun'
$ws.Range("K5").Value = 'Christina Ko'
$ws.Range("L5").Value = 'Susanna Tang'
$ws.Range("M5").Value = 'claudia lin'
$ws.Range("O5").Value = 'Jiwang Lee'
$ws.Range("C6").Value = 'Aaron Long'
$ws.Range("D6").Value = 'JJ Lee'
$ws.Range("E6").Value = 'Joann Jung'
$ws.Range("F6").Value = 'Jeffery Huang'
$ws.Range("G6").Value = 'Samuel Wen'
$ws.Range("K6").Value = 'Nathanael Wang'
$ws.Range("L6").Value = 'Logan Golia'
$ws.Range("M6").Value = 'Derek Liang '
$ws.Range("O6").Value = 'Israel Haile'
$ws.Range("E10").Value = 'Jiwang Lee'
$ws.Range("F10").Value = 'Justin Zhang'
$ws.Range("G10").Value = 'Logan Golia'
$ws.Range("K10").Value = 'Hannah Kim'
$ws.Range("L10").Value = 'Rachel Kim'
$ws.Range("E11").Value = 'Benjamin Kim'
$ws.Range("F11").Value = 'Daniel Kim '
$ws.Range("G11").Value = 'Jocelyn Youn'
$ws.Range("L11").Value = 'Claire Doh'
$ws.Range("N11").Value = 'Aaron Long'
$ws.Range("E12").Value = 'Grace Kwon'
$ws.Range("F12").Value = 'Daniel Song'
$ws.Range("G12").Value = 'Jeff Jiang'
$ws.Range("K12").Value = 'David Zhu'
$ws.Range("L12").Value = 'Joann Jung'
$ws.Range("N12").Value = 'Michael Dong '
$ws.Range("E13").Value = 'Claire Doh'
$ws.Range("F13").Value = 'Hannah Kim'
$ws.Range("G13").Value = 'Sam Ko'
$ws.Range("K13").Value = 'Grace Kwon'
$ws.Range("L13").Value = 'Daniel Kuo'
$ws.Range("C17").Value = 'Derek Liang '
$ws.Range("C18").Value = 'Taeho Choe'
$ws.Range("C19").Value = 'David Zhu'
$ws.Range("C20").Value = 'Israel Haile'
